# Apply weekly roll-forward update to the Fruta/Hortaliza "Membrillo" sheet.
# Row 2 takes former row 3's values, row 3 takes former row 4's values,
# and row 4 is repopulated with the former row 2's Date/Volume, keeping
# the other (N/O/P/R/S) values that row 4 already had.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (was: 45043 / Primera / 270 / 19000 / 20000 / 19500 / O'Higgins / 1083)
$ws.Range("D2").Value = 45086
$ws.Range("L2").Value = "Segunda"
$ws.Range("M2").Value = 270
$ws.Range("N2").Value = 20000
$ws.Range("O2").Value = 21000
$ws.Range("P2").Value = 20500
$ws.Range("R2").Value = "Provincia de Curicó"
$ws.Range("S2").Value = 1139

# --- Row 3 (was: 45086 / Segunda / 270 / 20000 / 21000 / 20500 / Curicó / 1139)
$ws.Range("D3").Value = 45107
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 200
$ws.Range("N3").Value = 19000
$ws.Range("O3").Value = 20000
$ws.Range("P3").Value = 19500
$ws.Range("R3").Value = "Región de O'Higgins"
$ws.Range("S3").Value = 1083

# --- Row 4 (was: 45107 / Primera / 200 / 19000 / 20000 / 19500 / O'Higgins / 1083)
$ws.Range("D4").Value = 45043
$ws.Range("M4").Value = 270
